# Weekly update: insert two new rows of data (new week) above the existing
# historical rows 37-44, shifting them down to rows 39-46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 37 and 38 (existing rows 37-44 shift to 39-46)
$ws.Range("A37:A38").EntireRow.Insert()

# --- New row 37 ---
$ws.Cells.Item(37, 1).Value = 1
$ws.Cells.Item(37, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(37, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(37, 4).Value = 44889
$ws.Cells.Item(37, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(37, 5).Value = 15
$ws.Cells.Item(37, 6).Value = 100112045
$ws.Cells.Item(37, 7).Value = "Zapallo"
$ws.Cells.Item(37, 8).Value = "Camote"
$ws.Cells.Item(37, 9).Value = "1a nueva(o)"
$ws.Cells.Item(37, 10).Value = 500
$ws.Cells.Item(37, 11).Value = 830
$ws.Cells.Item(37, 12).Value = 850
$ws.Cells.Item(37, 13).Value = 838
$ws.Cells.Item(37, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(37, 15).Value = "Perú"
$ws.Cells.Item(37, 16).Value = 838
$ws.Cells.Item(37, 17).Value = 1
$ws.Cells.Item(37, 18).Value = "Hortaliza"

# --- New row 38 ---
$ws.Cells.Item(38, 1).Value = 1
$ws.Cells.Item(38, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(38, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(38, 4).Value = 44889
$ws.Cells.Item(38, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(38, 5).Value = 15
$ws.Cells.Item(38, 6).Value = 100112045
$ws.Cells.Item(38, 7).Value = "Zapallo"
$ws.Cells.Item(38, 8).Value = "Camote"
$ws.Cells.Item(38, 9).Value = "2a nueva(o)"
$ws.Cells.Item(38, 10).Value = 400
$ws.Cells.Item(38, 11).Value = 830
$ws.Cells.Item(38, 12).Value = 850
$ws.Cells.Item(38, 13).Value = 842
$ws.Cells.Item(38, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(38, 15).Value = "Perú"
$ws.Cells.Item(38, 16).Value = 842
$ws.Cells.Item(38, 17).Value = 1
$ws.Cells.Item(38, 18).Value = "Hortaliza"
